# Adds the "Totenkopf / Firewall anschauen / ... / GESCHAFFT!" ending block
# to the escape-room script, replacing the single trailing empty paragraph
# with 14 new paragraphs (body text, Überschrift2 headings, a bulleted list
# item and a few tab-separated reference lines) taken from the target diff.

$d = $word.ActiveDocument
$wNs = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# `InsertXML` replaces the *content* of a range with literal OOXML, which lets us
# reproduce the exact run/proofErr/pPr structure from the diff. The document
# currently ends with a single empty paragraph (right before the sectPr) -
# that paragraph becomes the first new paragraph below; every subsequent
# paragraph is created with InsertParagraphAfter() and then filled the same way.
$target = $d.Paragraphs.Last.Range

$xml0 = '<w:p' + $wNs + '><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Der Totenkopf verschwindet. Du hast es geschafft. Doch dir fällt etwas ein: „Wir müssen irgendetwas tun, um die Hacker aus dem System zu werfen und das System </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>besser.aBSIchern</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>!“, rufst du. Du wendest dich wieder dem Kontrollrechner zu um dir den Status der Firewall anzuschauen.</w:t></w:r></w:p>'
$target.InsertXML($xml0)

$xml1 = '<w:p' + $wNs + '><w:pPr><w:pStyle w:val="berschrift2"/></w:pPr><w:r><w:t>Firewall anschauen</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml1)

$xml2 = '<w:p' + $wNs + '><w:r><w:t>„</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ohje</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>…jede Menge Lücken</w:t></w:r><w:r><w:t>! Wir müssen irgendetwas tun, um die Lücken zu schließen:</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml2)

$xml3 = '<w:p' + $wNs + '><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Firewall-ASCII wird mit Lücken angezeigt</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml3)

$xml4 = '<w:p' + $wNs + '><w:pPr><w:pStyle w:val="berschrift2"/></w:pPr><w:r><w:t>3. Umschauen Raum 6</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml4)

$xml5 = '<w:p' + $wNs + '><w:r><w:t>In einer Ecke des Kontrollpultes liegt ein Zettel. Vielleicht hilft dir dieser ja weiter.</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml5)

$xml6 = '<w:p' + $wNs + '><w:pPr><w:pStyle w:val="berschrift2"/></w:pPr><w:r><w:t>Zettel nehmen</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml6)

$xml7 = '<w:p' + $wNs + '><w:r><w:t xml:space="preserve">ORP.3.A7 </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>|</w:t></w:r><w:r><w:tab/><w:t>5</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml7)

$xml8 = '<w:p' + $wNs + '><w:r><w:t>APP.3.4</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>| 2.2</w:t></w:r><w:r><w:tab/><w:t>|</w:t></w:r><w:r><w:tab/><w:t>14</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml8)

$xml9 = '<w:p' + $wNs + '><w:r><w:t>ISMS.1.A11</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>|</w:t></w:r><w:r><w:tab/><w:t>4</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml9)

$xml10 = '<w:p' + $wNs + '><w:r><w:t>OPS.1.1.4.A14</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>|</w:t></w:r><w:r><w:tab/><w:t>5</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml10)

$xml11 = '<w:p' + $wNs + '><w:r><w:t>SYS.2.1.A1</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>|</w:t></w:r><w:r><w:tab/><w:t>5</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml11)

$xml12 = '<w:p' + $wNs + '><w:pPr><w:pStyle w:val="berschrift2"/></w:pPr><w:r><w:t>Satz eingeben/Firewall schließen</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml12)

$xml13 = '<w:p' + $wNs + '><w:r><w:t xml:space="preserve">Super. Du hast die Angreifer ausgesperrt und die Lücken in der Firewall geschlossen. Ab hier übernimmt der Kraftwerkchef. In letzter Sekunde fährt er über </w:t></w:r><w:r><w:t>den Kontrollrechner die Pumpen des Kühlsystems wieder hoch. Das rote Notlicht erlischt und das Warnsignal aus dem Maschinenraum ist auch nicht mehr zu hören. GESCHAFFT!</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last.Range
$target.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml13)

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
